$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5-6 down to 6-7
$ws.Rows(5).Insert()

# Fill the new row 5 with the new daily price entry
$ws.Range("A5").Value = 8
$ws.Range("B5").Value = "Terminal La Palmera de La Serena"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44438
$ws.Range("D5").NumberFormat = $ws.Range("D6").NumberFormat
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 100112013
$ws.Range("G5").Value = "Alcachofa"
$ws.Range("H5").Value = "Española"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 11500
$ws.Range("N5").Value = "`$/caja 30 unidades"
$ws.Range("O5").Value = "Provincia del Elquí"
$ws.Range("P5").Value = 383
$ws.Range("Q5").Value = 30
$ws.Range("R5").Value = "Hortaliza"
